# update terbaru 8 juni 2023
$wb = $excel.ActiveWorkbook

# "Global" sheet holds the cover subtitle text with the build number
$ws = $wb.Worksheets.Item("Global")

# Update the build number in the cover subtitle from (561) to (568)
$ws.Range("B2").Value = "BSI Super Apps - App Version 1.0.2 (568) Rebrand OCP QA"

# Move the active selection on the sheet from B3 to B2
$ws.Activate()
$ws.Range("B2").Select()
